$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E25) from "S.GISH" to "fullRNASEQ"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
}

# Scroll/selection state: top-left cell A19, selection E24:E25 with active cell E24
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E24:E25").Select()

# Enable iterative calculation with a max change (iterateDelta) of 1E-4
$wb.Application.Iteration = $true
$wb.Application.MaxChange = 0.0001
